$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038927444142529
$ws.Range("D2").Value = 1.047616587205179
$ws.Range("E2").Value = 1.037524789852911
$ws.Range("F2").Value = 1.056188527193665
$ws.Range("I2").Value = 1.042043491777227
$ws.Range("J2").Value = 1.044022227579138
$ws.Range("K2").Value = 1.050378726037938
$ws.Range("L2").Value = 1.040315426889858
$ws.Range("M2").Value = 1.058926923289173
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039787584377955
$ws.Range("D3").Value = 1.048301733414593
$ws.Range("E3").Value = 1.03825420163419
$ws.Range("F3").Value = 1.057004621183756
$ws.Range("I3").Value = 1.042256228456884
$ws.Range("J3").Value = 1.04452781047057
$ws.Range("K3").Value = 1.05087608095883
$ws.Range("L3").Value = 1.040854877054477
$ws.Range("M3").Value = 1.059556600591692
$ws.Range("N3").Value = 1.005712725503983
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040344838241483
$ws.Range("D4").Value = 1.048745672616951
$ws.Range("E4").Value = 1.038727138843579
$ws.Range("F4").Value = 1.057533635725298
$ws.Range("I4").Value = 1.042393072671594
$ws.Range("J4").Value = 1.044854980474361
$ws.Range("K4").Value = 1.051197826453084
$ws.Range("L4").Value = 1.041204226393943
$ws.Range("M4").Value = 1.059964336641698
$ws.Range("N4").Value = 1.005712725503983
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040579270305735
$ws.Range("D5").Value = 1.048932447431525
$ws.Range("E5").Value = 1.03892618944282
$ws.Range("F5").Value = 1.057756258273463
$ws.Range("I5").Value = 1.042450407178453
$ws.Range("J5").Value = 1.044992527229103
$ws.Range("K5").Value = 1.05133306863562
$ws.Range("L5").Value = 1.041351160674433
$ws.Range("M5").Value = 1.06013581737374
$ws.Range("N5").Value = 1.005712725503983
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040618641965935
$ws.Range("D6").Value = 1.048963816065455
$ws.Range("E6").Value = 1.038959624220237
$ws.Range("F6").Value = 1.057793650707586
$ws.Range("I6").Value = 1.042460022449987
$ws.Range("J6").Value = 1.045015622160384
$ws.Range("K6").Value = 1.051355775203685
$ws.Range("L6").Value = 1.041375835523683
$ws.Range("M6").Value = 1.060164613713851
$ws.Range("N6").Value = 1.005712725503983
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040347970097522
$ws.Range("D7").Value = 1.048748167752391
$ws.Range("E7").Value = 1.038729797674598
$ws.Range("F7").Value = 1.057536609537641
$ws.Range("I7").Value = 1.042393839544426
$ws.Range("J7").Value = 1.044856818364805
$ws.Range("K7").Value = 1.051199633644524
$ws.Range("L7").Value = 1.041206189472233
$ws.Range("M7").Value = 1.059966627708473
$ws.Range("N7").Value = 1.005712725503983
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039217989639332
$ws.Range("D8").Value = 1.047848009488851
$ws.Range("E8").Value = 1.037771098331759
$ws.Range("F8").Value = 1.056464132710306
$ws.Range("I8").Value = 1.042115554542539
$ws.Range("J8").Value = 1.044193085933899
$ws.Range("K8").Value = 1.050546824478413
$ws.Range("L8").Value = 1.040497675932355
$ws.Range("M8").Value = 1.059139664117747
$ws.Range("N8").Value = 1.005712725503983
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.037232141463659
$ws.Range("D9").Value = 1.046266518410763
$ws.Range("E9").Value = 1.036089171781897
$ws.Range("F9").Value = 1.054581630431147
$ws.Range("I9").Value = 1.041619004688499
$ws.Range("J9").Value = 1.043023745370419
$ws.Range("K9").Value = 1.049395965223945
$ws.Range("L9").Value = 1.039251458697217
$ws.Range("M9").Value = 1.057684760637539
$ws.Range("N9").Value = 1.005712725503983
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035911915938127
$ws.Range("D10").Value = 1.045215458859367
$ws.Range("E10").Value = 1.034972986069888
$ws.Range("F10").Value = 1.053331677643743
$ws.Range("I10").Value = 1.041283862718402
$ws.Range("J10").Value = 1.042244423016686
$ws.Range("K10").Value = 1.048628451605641
$ws.Range("L10").Value = 1.038422258061563
$ws.Range("M10").Value = 1.056716473988286
$ws.Range("N10").Value = 1.005712725503983
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035341134963285
$ws.Range("D11").Value = 1.044761135949633
$ws.Range("E11").Value = 1.034490897304001
$ws.Range("F11").Value = 1.052791655963897
$ws.Range("I11").Value = 1.041137777226266
$ws.Range("J11").Value = 1.041907040342505
$ws.Range("K11").Value = 1.048296061516011
$ws.Range("L11").Value = 1.038063604384843
$ws.Range("M11").Value = 1.056297606258741
$ws.Range("N11").Value = 1.005712725503983
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035129256117221
$ws.Range("D12").Value = 1.044592501117122
$ws.Range("E12").Value = 1.034312014254794
$ws.Range("F12").Value = 1.052591252536391
$ws.Range("I12").Value = 1.04108337004742
$ws.Range("J12").Value = 1.041781733058563
$ws.Range("K12").Value = 1.048172590630501
$ws.Range("L12").Value = 1.037930445279823
$ws.Range("M12").Value = 1.056142082857525
$ws.Range("N12").Value = 1.005712725503983
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035174698731243
$ws.Range("D13").Value = 1.044628668359924
$ws.Range("E13").Value = 1.034350376821383
$ws.Range("F13").Value = 1.052634231363757
$ws.Range("I13").Value = 1.041095047100193
$ws.Range("J13").Value = 1.041808611347836
$ws.Range("K13").Value = 1.048199075822803
$ws.Range("L13").Value = 1.037959005584361
$ws.Range("M13").Value = 1.056175440290798
$ws.Range("N13").Value = 1.005712725503983
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035323618236243
$ws.Range("D14").Value = 1.044747194064455
$ws.Range("E14").Value = 1.034476106966637
$ws.Range("F14").Value = 1.052775086775541
$ws.Range("I14").Value = 1.041133282855881
$ws.Range("J14").Value = 1.041896682157937
$ws.Range("K14").Value = 1.04828585549698
$ws.Range("L14").Value = 1.038052596162509
$ws.Range("M14").Value = 1.056284749357026
$ws.Range("N14").Value = 1.005712725503983
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035415390344962
$ws.Range("D15").Value = 1.044820237740815
$ws.Range("E15").Value = 1.03455359818483
$ws.Range("F15").Value = 1.05286189697259
$ws.Range("I15").Value = 1.041156822038105
$ws.Range("J15").Value = 1.041950947074821
$ws.Range("K15").Value = 1.048339322509792
$ws.Range("L15").Value = 1.03811026851013
$ws.Range("M15").Value = 1.056352106640029
$ws.Range("N15").Value = 1.005712725503983
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035949815582413
$ws.Range("D16").Value = 1.045245627633101
$ws.Range("E16").Value = 1.035005006722287
$ws.Range("F16").Value = 1.053367542893764
$ws.Range("I16").Value = 1.04129353762407
$ws.Range("J16").Value = 1.04226681555844
$ws.Range("K16").Value = 1.048650510268661
$ws.Range("L16").Value = 1.03844606918025
$ws.Range("M16").Value = 1.056744281567744
$ws.Range("N16").Value = 1.005712725503983
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036285284588226
$ws.Range("D17").Value = 1.045512677012953
$ws.Range("E17").Value = 1.03528849317637
$ws.Range("F17").Value = 1.053685048195893
$ws.Range("I17").Value = 1.041379037371463
$ws.Range("J17").Value = 1.042464970909185
$ws.Range("K17").Value = 1.048845697258794
$ws.Range("L17").Value = 1.038656814985568
$ws.Range("M17").Value = 1.056990392581254
$ws.Range("N17").Value = 1.005712725503983
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036481043303234
$ws.Range("D18").Value = 1.045668518771379
$ws.Range("E18").Value = 1.03545396429457
$ws.Range("F18").Value = 1.053870360873913
$ws.Range("I18").Value = 1.041428814595642
$ws.Range("J18").Value = 1.042580558148325
$ws.Range("K18").Value = 1.0489595414559
$ws.Range("L18").Value = 1.038779777531209
$ws.Range("M18").Value = 1.05713398411189
$ws.Range("N18").Value = 1.005712725503983
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036547806389197
$ws.Range("D19").Value = 1.045721669682723
$ws.Range("E19").Value = 1.035510405673685
$ws.Range("F19").Value = 1.053933567528937
$ws.Range("I19").Value = 1.041445771489859
$ws.Range("J19").Value = 1.042619971485941
$ws.Range("K19").Value = 1.048998358483767
$ws.Range("L19").Value = 1.038821710998734
$ws.Range("M19").Value = 1.057182951703015
$ws.Range("N19").Value = 1.005712725503983
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036249283102182
$ws.Range("D20").Value = 1.045484017229096
$ws.Range("E20").Value = 1.035258065516693
$ws.Range("F20").Value = 1.053650970744405
$ws.Range("I20").Value = 1.041369873706232
$ws.Range("J20").Value = 1.042443710025929
$ws.Range("K20").Value = 1.04882475605686
$ws.Range("L20").Value = 1.038634200005979
$ws.Range("M20").Value = 1.056963983138598
$ws.Range("N20").Value = 1.005712725503983
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035279761424626
$ws.Range("D21").Value = 1.044712287845965
$ws.Range("E21").Value = 1.034439077418267
$ws.Range("F21").Value = 1.052733603253996
$ws.Range("I21").Value = 1.041122027362495
$ws.Range("J21").Value = 1.041870747167146
$ws.Range("K21").Value = 1.048260301211868
$ws.Range("L21").Value = 1.038025034379259
$ws.Range("M21").Value = 1.056252558817069
$ws.Range("N21").Value = 1.005712725503983
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034670964343533
$ws.Range("D22").Value = 1.044227771484142
$ws.Range("E22").Value = 1.033925225204375
$ws.Range("F22").Value = 1.052157887309715
$ws.Range("I22").Value = 1.040965360728593
$ws.Range("J22").Value = 1.041510570795771
$ws.Range("K22").Value = 1.047905369650113
$ws.Range("L22").Value = 1.037642380461291
$ws.Range("M22").Value = 1.055805621806513
$ws.Range("N22").Value = 1.005712725503983
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034993624754437
$ws.Range("D23").Value = 1.044484555767249
$ws.Range("E23").Value = 1.034197525179495
$ws.Range("F23").Value = 1.052462983202846
$ws.Range("I23").Value = 1.041048491655581
$ws.Range("J23").Value = 1.04170150021728
$ws.Range("K23").Value = 1.048093528585626
$ws.Range("L23").Value = 1.03784519864674
$ws.Range("M23").Value = 1.056042516560258
$ws.Range("N23").Value = 1.005712725503983
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036265550360625
$ws.Range("D24").Value = 1.045496967116028
$ws.Range("E24").Value = 1.035271814101175
$ws.Range("F24").Value = 1.053666368515986
$ws.Range("I24").Value = 1.041374014660574
$ws.Range("J24").Value = 1.04245331688375
$ws.Range("K24").Value = 1.04883421850045
$ws.Range("L24").Value = 1.03864441862431
$ws.Range("M24").Value = 1.056975916308218
$ws.Range("N24").Value = 1.005712725503983
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037744890366884
$ws.Range("D25").Value = 1.046674803537988
$ws.Range("E25").Value = 1.036523099271999
$ws.Range("F25").Value = 1.055067420901252
$ws.Range("I25").Value = 1.041748102131622
$ws.Range("J25").Value = 1.043326011559769
$ws.Range("K25").Value = 1.049693544073043
$ws.Range("L25").Value = 1.039573357792826
$ws.Range("M25").Value = 1.058060604336421
$ws.Range("N25").Value = 1.005712725503983

Write-Output "applied vm_pu update"
